$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 381 (shifts existing 381:500 down to 383:502,
# and Excel's dimension auto-extends to A1:T502).
$ws.Rows("381:382").Insert()

# New row 381: Naranja / Lane Late / Primera, week of 2022-12-29 (serial 44924)
$ws.Cells.Item(381, 1).Value2 = 2
$ws.Cells.Item(381, 2).Value2 = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(381, 3).Value2 = "Coquimbo"
$ws.Cells.Item(381, 4).Value2 = 44924
$ws.Cells.Item(381, 5).Value2 = 4
$ws.Cells.Item(381, 6).Value2 = "Fruta"
$ws.Cells.Item(381, 7).Value2 = 100102
$ws.Cells.Item(381, 8).Value2 = "Cítricos"
$ws.Cells.Item(381, 9).Value2 = 100102005
$ws.Cells.Item(381, 10).Value2 = "Naranja"
$ws.Cells.Item(381, 11).Value2 = "Lane Late"
$ws.Cells.Item(381, 12).Value2 = "Primera"
$ws.Cells.Item(381, 13).Value2 = 20
$ws.Cells.Item(381, 14).Value2 = 200000
$ws.Cells.Item(381, 15).Value2 = 210000
$ws.Cells.Item(381, 16).Value2 = 205000
$ws.Cells.Item(381, 17).Value2 = "$/bins (400 kilos)"
$ws.Cells.Item(381, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(381, 19).Value2 = 512
$ws.Cells.Item(381, 20).Value2 = 400

# New row 382: Naranja / Lane Late / Segunda, same week
$ws.Cells.Item(382, 1).Value2 = 2
$ws.Cells.Item(382, 2).Value2 = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(382, 3).Value2 = "Coquimbo"
$ws.Cells.Item(382, 4).Value2 = 44924
$ws.Cells.Item(382, 5).Value2 = 4
$ws.Cells.Item(382, 6).Value2 = "Fruta"
$ws.Cells.Item(382, 7).Value2 = 100102
$ws.Cells.Item(382, 8).Value2 = "Cítricos"
$ws.Cells.Item(382, 9).Value2 = 100102005
$ws.Cells.Item(382, 10).Value2 = "Naranja"
$ws.Cells.Item(382, 11).Value2 = "Lane Late"
$ws.Cells.Item(382, 12).Value2 = "Segunda"
$ws.Cells.Item(382, 13).Value2 = 16
$ws.Cells.Item(382, 14).Value2 = 160000
$ws.Cells.Item(382, 15).Value2 = 170000
$ws.Cells.Item(382, 16).Value2 = 165000
$ws.Cells.Item(382, 17).Value2 = "$/bins (400 kilos)"
$ws.Cells.Item(382, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(382, 19).Value2 = 412
$ws.Cells.Item(382, 20).Value2 = 400
